# Generate Report for Archive
# - Update the handoff status text from "Ready for handoff" to "In Translation"
#   (this shows up on the Overview sheet's zh-cn/de-de status columns and on
#   each language sheet's own Status column).
# - Shrink the now-narrower Status columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 13.4101845877511

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
